$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.139.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.313.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.92%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.78%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.667.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.284.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.806"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.051.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("E19").Value = "  +6.86%  "
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +14.46%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.97%  "
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.994.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("E45").Value = "  +5.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.49%  "
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.534.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
